# 案件情報.xlsx — append new scraped listing at the top of the detail rows
# (row 9) and refresh the "取得日時" (fetched-at) timestamp on every data
# row to the new scrape time, 2025-10-21 01:47:06.
#
# The sheet "ランサーズ" is a simple table: row 1 is the header, rows 2..25
# hold the scraped listings (newest/highest score first). This edit:
#   1. Inserts a brand-new row at position 9 (pushing the former rows
#      9..25 down to 10..26), matching the new dimension A1:H26.
#   2. Fills that new row 9 with the freshly scraped listing.
#   3. Stamps column A (取得日時) on every data row (2..26) with the new
#      scrape timestamp, since the whole scrape batch shares one run time.
#   4. Rebuilds the F-column hyperlinks so every URL cell again carries a
#      live hyperlink (Rows.Insert()/Range.Hyperlinks in this host mangles
#      per-cell hyperlinks, so the safest path is to drop the sheet's
#      hyperlink collection and re-add one per row from the cell text).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = '2025-10-21 01:47:06'

# 1. Push rows 9..25 down to 10..26, opening up a fresh row 9.
$ws.Rows.Item(9).Insert()

# 2. Populate the newly opened row 9 with the new listing.
$ws.Range('A9').Value2 = $newTimestamp
$ws.Range('B9').Value2 = '【急募】ebayAPIを活用したShippingポリシー設定の専門家募集'
$ws.Range('C9').Value2 = 'システム開発'
$ws.Range('D9').Value2 = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range('E9').Value2 = '期限情報なし'
$ws.Range('F9').Value2 = 'https://www.lancers.jp/work/detail/5415908'
$ws.Range('G9').Value2 = 183
$ws.Range('H9').Value2 = '🔥API'

# 3. Refresh the fetched-at timestamp for every data row (2..26) — the
#    rows that existed before the insert as well as the new row 9.
for ($r = 2; $r -le 26; $r++) {
    $ws.Range("A$r").Value2 = $newTimestamp
}

# 4. Rebuild hyperlinks on F2:F26 from the cell text, in row order, so
#    relationship ids line up 1:1 with the visible rows (rId1 -> F2, ...,
#    rId25 -> F26) — mirrors how the previous hyperlink set was laid out.
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le 26; $r++) {
    $cell = $ws.Range("F$r")
    $ws.Hyperlinks.Add($cell, $cell.Value2)
}
